$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-05-28"

# Update the header label in I1 to reflect the new "through" date
$ws.Range("I1").Value = "2022 (through 05-28)"

# Update the June (row 6) "Total" column value
$ws.Range("I6").Value = 101

# Update the yearly Total (row 14) "Total" column value
$ws.Range("I14").Value = 652
